$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = "13-11-2025"
$newPrice = "The price of gold in India today is ₹12,780 per gram for 24 karat gold, ₹11,715 per gram for 22 karat gold and ₹9,585 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A58").Value = $newDate
$ws.Range("B58").Value = $newPrice

# Match formatting of the prior row (border on both, wrap text on column B only)
$ws.Range("A58:B58").Borders.LineStyle = 1
$ws.Range("B58").WrapText = $true

$wb.Save()
